# edit.ps1 - apply the changes described by the commit diff to $ppt.ActivePresentation
#
# Summary of changes:
#  1. Cached "datetimeFigureOut" footer field text 2017/3/2 -> 2017/6/16
#     (slide master, notes master, and every slide layout).
#  2. Merge the two runs "define " + "bean & AOP" into a single run
#     "define bean & AOP" (slides 9, 10 and 11 - the HandlerException
#     description diagrams).
#  3. Slide 12 (exception-level legend): move the "error" right-arrow up,
#     and delete the now redundant "warn" / "Info" right-arrows.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: find a shape in a Shapes collection by its stable COM Id.
# ---------------------------------------------------------------------
function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# Helper: find the single shape in a Shapes collection whose text equals
# $text exactly.
function Get-ShapeByText($shapes, $text) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq $text) {
            return $sh
        }
    }
    return $null
}

# Helper: force-set text - going through a throwaway value first makes sure
# the run is actually rewritten even in the (rare) case the target text is
# already equal to the current text.
function Set-ShapeText($shape, $newText) {
    $shape.TextFrame.TextRange.Text = "~~tmp~~"
    $shape.TextFrame.TextRange.Text = $newText
}

# ---------------------------------------------------------------------
# 1. Footer date placeholder: 2017/3/2 -> 2017/6/16
#    (slide master + notes master + all custom (slide) layouts)
# ---------------------------------------------------------------------
$dateContainers = New-Object System.Collections.ArrayList
[void]$dateContainers.Add($p.SlideMaster.Shapes)
[void]$dateContainers.Add($p.NotesMaster.Shapes)
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    [void]$dateContainers.Add($p.SlideMaster.CustomLayouts.Item($li).Shapes)
}

foreach ($shapes in $dateContainers) {
    $dateShape = Get-ShapeByText $shapes "2017/3/2"
    if ($dateShape -ne $null) {
        Set-ShapeText $dateShape "2017/6/16"
    }
}

# ---------------------------------------------------------------------
# 2. "define " + "bean & AOP" -> single run "define bean & AOP"
#    on slides 9, 10 and 11.
# ---------------------------------------------------------------------
foreach ($slideIdx in 9, 10, 11) {
    $s = $p.Slides.Item($slideIdx)
    $sh = Get-ShapeByText $s.Shapes "define bean & AOP"
    if ($sh -ne $null) {
        Set-ShapeText $sh "define bean & AOP"
    }
}

# ---------------------------------------------------------------------
# 3. Slide 12: reposition the "error" arrow and drop "warn" / "Info".
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)

$errorArrow = Get-ShapeById $s12.Shapes 76
if ($errorArrow -ne $null) {
    $errorArrow.Top = 25.95292
}

$warnArrow = Get-ShapeById $s12.Shapes 80
if ($warnArrow -ne $null) {
    $warnArrow.Delete()
}

$infoArrow = Get-ShapeById $s12.Shapes 81
if ($infoArrow -ne $null) {
    $infoArrow.Delete()
}
